$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# Row 2
$ws.Range("A2").Value = 13440
$ws.Range("B2").Value = "Vicente Silveira"
$ws.Range("C2").Value = "Marketing"
$ws.Range("D2").Value = "Viagem de negócios"
$ws.Range("E2").Value = 7
$ws.Range("F2").Value = 45082
$ws.Range("G2").Value = 6135.95

# Row 3
$ws.Range("A3").Value = 30657
$ws.Range("B3").Value = "Carolina da Mota"
$ws.Range("C3").Value = "P&D"
$ws.Range("D3").Value = "Consulta médica"
$ws.Range("E3").Value = 8
$ws.Range("F3").Value = 45090
$ws.Range("G3").Value = 7718.88

# Row 4
$ws.Range("A4").Value = 24439
$ws.Range("B4").Value = "Vitória da Cunha"
$ws.Range("C4").Value = "Vendas"
$ws.Range("E4").Value = 1
$ws.Range("F4").Value = 45086
$ws.Range("G4").Value = 10282.69

# Row 5
$ws.Range("A5").Value = 3829
$ws.Range("B5").Value = "Manuela da Mota"
$ws.Range("C5").Value = "TI"
$ws.Range("E5").Value = 1
$ws.Range("F5").Value = 45096
$ws.Range("G5").Value = 5408.03

# Row 6
$ws.Range("A6").Value = 71013
$ws.Range("B6").Value = "Dra. Isabel Caldeira"
$ws.Range("C6").Value = "Vendas"
$ws.Range("D6").Value = "Outros"
$ws.Range("E6").Value = 8
$ws.Range("G6").Value = 4685.06

# Row 7
$ws.Range("A7").Value = 8530
$ws.Range("B7").Value = "Dr. Davi Luiz Campos"
$ws.Range("C7").Value = "TI"
$ws.Range("F7").Value = 45104
$ws.Range("G7").Value = 7777.78

# Row 8
$ws.Range("A8").Value = 42489
$ws.Range("B8").Value = "Lavínia Duarte"
$ws.Range("E8").Value = 2
$ws.Range("F8").Value = 45082
$ws.Range("G8").Value = 6447.23

# Row 9
$ws.Range("A9").Value = 40095
$ws.Range("B9").Value = "Ana Lívia Melo"
$ws.Range("D9").Value = "Problemas pessoais"
$ws.Range("E9").Value = 2
$ws.Range("F9").Value = 45084
$ws.Range("G9").Value = 9692.73

# Row 10
$ws.Range("A10").Value = 8472
$ws.Range("B10").Value = "Igor Moura"
$ws.Range("C10").Value = "Recursos Humanos"
$ws.Range("D10").Value = "Viagem de negócios"
$ws.Range("F10").Value = 45102
$ws.Range("G10").Value = 5388.62

# Row 11
$ws.Range("A11").Value = 14979
$ws.Range("B11").Value = "Sra. Milena Cardoso"
$ws.Range("C11").Value = "Operações"
$ws.Range("D11").Value = "Viagem de negócios"
$ws.Range("E11").Value = 2
$ws.Range("F11").Value = 45099
$ws.Range("G11").Value = 5549.26
